$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Modified the born (RelivePos) position of the city "villageScene" (row 2) in the Scene config.
$ws.Range("E2").Value = "20,0,-137"
